$wb = $excel.ActiveWorkbook

# Rename the worksheets (Russian -> English). Renaming via the Excel
# object model automatically updates any formulas elsewhere in the
# workbook that reference these sheets by name.
$wb.Worksheets.Item("Сводный").Name = "Consolidated budget"
$wb.Worksheets.Item("Продажи").Name = "Sales"
$wb.Worksheets.Item("Производство").Name = "Production"
$wb.Worksheets.Item("Логистика").Name = "Logistics"
$wb.Worksheets.Item("Прочее").Name = "Misc"
